$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the very top; everything else shifts down one row.
$ws.Rows.Item(1).Insert()

# Populate the new header row (row 1) with the column titles.
$ws.Range("A1").Value = "CNE"
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"
$ws.Range("D1").Value = "DateofBirth"
$ws.Range("E1").Value = "ClasseName"
$ws.Range("F1").Value = "Phone"
$ws.Range("G1").Value = "Email"

# Update the first student id (now in row 2); the formulas below it
# (A3:A11 = previous + 1) recompute automatically.
$ws.Range("A2").Value = 17000041

# Match the cursor position left behind in the saved file.
$ws.Range("I9").Select()
